$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Jonas Valanciunas", "C", "Washington Wizards"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
